$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header change in D1
$ws.Range("D1").Value = "Unnamed: 3"

# Updated numeric values in column D, rows 2-15
$ws.Range("D2").Value = 525.11
$ws.Range("D3").Value = 177.46
$ws.Range("D4").Value = 197.56
$ws.Range("D5").Value = 243.31
$ws.Range("D6").Value = 249.3
$ws.Range("D7").Value = 255.72
$ws.Range("D8").Value = 275.39
$ws.Range("D9").Value = 277.1
$ws.Range("D10").Value = 285.65
$ws.Range("D11").Value = 292.92
$ws.Range("D12").Value = 300.62
$ws.Range("D13").Value = 307.46
$ws.Range("D14").Value = 313.87
$ws.Range("D15").Value = 313.87
